$d = $word.ActiveDocument

$d.Content.Find.Execute("56×42=2352", $true, $false, $false, $false, $false, $true, 1, $false, "37×72=2664", 2)
$d.Content.Find.Execute("68×18=1224", $true, $false, $false, $false, $false, $true, 1, $false, "92×79=7268", 2)
$d.Content.Find.Execute("23×42=966", $true, $false, $false, $false, $false, $true, 1, $false, "79×60=4740", 2)
$d.Content.Find.Execute("65×55=3575", $true, $false, $false, $false, $false, $true, 1, $false, "94×76=7144", 2)
$d.Content.Find.Execute("33×98=3234", $true, $false, $false, $false, $false, $true, 1, $false, "49×92=4508", 2)
$d.Content.Find.Execute("32×19=608", $true, $false, $false, $false, $false, $true, 1, $false, "63×68=4284", 2)
$d.Content.Find.Execute("81×24=1944", $true, $false, $false, $false, $false, $true, 1, $false, "33×32=1056", 2)
$d.Content.Find.Execute("57×63=3591", $true, $false, $false, $false, $false, $true, 1, $false, "94×94=8836", 2)
$d.Content.Find.Execute("84×74=6216", $true, $false, $false, $false, $false, $true, 1, $false, "95×77=7315", 2)
$d.Content.Find.Execute("41×72=2952", $true, $false, $false, $false, $false, $true, 1, $false, "33×57=1881", 2)
$d.Content.Find.Execute("33×20=660", $true, $false, $false, $false, $false, $true, 1, $false, "76×87=6612", 2)
$d.Content.Find.Execute("20×67=1340", $true, $false, $false, $false, $false, $true, 1, $false, "29×38=1102", 2)
$d.Content.Find.Execute("16×98=1568", $true, $false, $false, $false, $false, $true, 1, $false, "87×86=7482", 2)
$d.Content.Find.Execute("26×13=338", $true, $false, $false, $false, $false, $true, 1, $false, "90×81=7290", 2)
$d.Content.Find.Execute("47×64=3008", $true, $false, $false, $false, $false, $true, 1, $false, "72×72=5184", 2)
$d.Content.Find.Execute("45×16=720", $true, $false, $false, $false, $false, $true, 1, $false, "58×89=5162", 2)
$d.Content.Find.Execute("77×24=1848", $true, $false, $false, $false, $false, $true, 1, $false, "69×36=2484", 2)
$d.Content.Find.Execute("88×33=2904", $true, $false, $false, $false, $false, $true, 1, $false, "22×76=1672", 2)
$d.Content.Find.Execute("32×58=1856", $true, $false, $false, $false, $false, $true, 1, $false, "28×40=1120", 2)
$d.Content.Find.Execute("55×61=3355", $true, $false, $false, $false, $false, $true, 1, $false, "83×78=6474", 2)
$d.Content.Find.Execute("93×87=8091", $true, $false, $false, $false, $false, $true, 1, $false, "22×96=2112", 2)
$d.Content.Find.Execute("87×69=6003", $true, $false, $false, $false, $false, $true, 1, $false, "67×50=3350", 2)
$d.Content.Find.Execute("25×22=550", $true, $false, $false, $false, $false, $true, 1, $false, "96×84=8064", 2)
$d.Content.Find.Execute("53×26=1378", $true, $false, $false, $false, $false, $true, 1, $false, "33×61=2013", 2)
$d.Content.Find.Execute("14×22=308", $true, $false, $false, $false, $false, $true, 1, $false, "31×90=2790", 2)
